$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.964516129032258
$ws.Range("C2").Value = 0.8612903225806452

$ws.Range("B3").Value = 0.9612903225806452
$ws.Range("C3").Value = 0.8774193548387097

$ws.Range("B4").Value = 0.9709677419354839
$ws.Range("C4").Value = 0.8580645161290322

$ws.Range("B5").Value = 0.964516129032258
$ws.Range("C5").Value = 0.8806451612903226

$ws.Range("B6").Value = 0.964516129032258
$ws.Range("C6").Value = 0.8645161290322581
